# Update db parameters for freq reg demo
$wb = $excel.ActiveWorkbook

# Rename TGOV1N sheet to TGOV1DB
$ws = $wb.Worksheets.Item("TGOV1N")
$ws.Name = "TGOV1DB"

# Add new dbL / dbU columns (O, P) with header + values for the 10 data rows
$ws.Range("O1").Value = "dbL"
$ws.Range("P1").Value = "dbU"

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 15).Value = -0.0006
    $ws.Cells.Item($r, 16).Value = 0.0006
}

# Select the active cell / top-left per the recorded view state
$ws.Activate()
$ws.Range("N17").Select()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
